$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date rows are text strings like "2019-05-01" (not real dates), so we
# temporarily force a Text number format before assigning the value (this
# stops Excel's autodetect from turning the string into a date serial),
# then reset the style back to Normal so no stray style index is left on
# the cell (matches the source file, where these cells carry no "s" attr).

$ws.Range("C18:C22").NumberFormat = "@"

$ws.Range("C18").Value = "2019-05-01"
$ws.Range("E18").Value = 403

$ws.Range("C19").Value = "2019-05-11"
$ws.Range("E19").Value = 376

$ws.Range("C20").Value = "2019-05-21"

$ws.Range("C21").Value = "2019-05-31"

$ws.Range("C22").Value = "2019-06-01"
$ws.Range("E22").Value = 542

$ws.Range("C18:C22").Style = "Normal"
